$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 49.285713
$ws.Range("I38").Value = 49.285713
$ws.Range("K38").Value = 147.857139
$ws.Range("M38").Value = 224.142861

$ws.Range("H39").Value = 737.875
$ws.Range("I39").Value = 517
$ws.Range("J39").Value = 1106
$ws.Range("K39").Value = 1551
$ws.Range("L39").Value = 3318
$ws.Range("M39").Value = -1255
$ws.Range("N39").Value = -3910

$ws.Range("H74").Value = 4466.6665
$ws.Range("I74").Value = 4360
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 4360
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -3424
$ws.Range("N74").Value = -6872

$ws.Range("H77").Value = 4466.6665
$ws.Range("I77").Value = 4360
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 21800
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -17120
$ws.Range("N77").Value = -34360

$ws.Range("H113").Value = 1675.375
$ws.Range("J113").Value = 1685.0769
$ws.Range("L113").Value = 1685.0769
$ws.Range("N113").Value = -8193.0769

$ws.Range("H137").Value = 27779994
$ws.Range("I137").Value = 1400
$ws.Range("J137").Value = 125005070
$ws.Range("K137").Value = 4200
$ws.Range("L137").Value = 375015210
$ws.Range("M137").Value = -1650
$ws.Range("N137").Value = -375020310

$ws.Range("H141").Value = 956.5263
$ws.Range("I141").Value = 851.41174
$ws.Range("J141").Value = 1850
$ws.Range("K141").Value = 2554.23522
$ws.Range("L141").Value = 5550
$ws.Range("M141").Value = 2625.76478
$ws.Range("N141").Value = -15910

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2821.8333
$ws.Range("I61").Value = 2245
$ws.Range("J61").Value = 4552.3335
$ws.Range("K61").Value = 2245
$ws.Range("L61").Value = 4552.3335
$ws.Range("M61").Value = -2033
$ws.Range("N61").Value = -4976.3335

$ws.Range("H122").Value = 1198.2354
$ws.Range("I122").Value = 1131
$ws.Range("J122").Value = 1359.6
$ws.Range("K122").Value = 3393
$ws.Range("L122").Value = 4078.8
$ws.Range("M122").Value = -943
$ws.Range("N122").Value = -8978.799999999999

$ws.Range("H132").Value = 1735.1765
$ws.Range("I132").Value = 1586.138
$ws.Range("J132").Value = 2599.6
$ws.Range("K132").Value = 4758.414
$ws.Range("L132").Value = 7798.799999999999
$ws.Range("M132").Value = -2228.414
$ws.Range("N132").Value = -12858.8

$ws.Range("H136").Value = 2821.8333
$ws.Range("I136").Value = 2245
$ws.Range("J136").Value = 4552.3335
$ws.Range("K136").Value = 6735
$ws.Range("L136").Value = 13657.0005
$ws.Range("M136").Value = -4185
$ws.Range("N136").Value = -18757.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 50002000
$ws.Range("I7").Value = 50002000
$ws.Range("K7").Value = 50002000
$ws.Range("M7").Value = -50001887

$ws.Range("H20").Value = 1704.7391
$ws.Range("I20").Value = 1606.4117
$ws.Range("J20").Value = 1983.3334
$ws.Range("K20").Value = 1606.4117
$ws.Range("L20").Value = 1983.3334
$ws.Range("M20").Value = -1359.4117
$ws.Range("N20").Value = -2477.3334

$ws.Range("H99").Value = 2017.2
$ws.Range("I99").Value = 890
$ws.Range("K99").Value = 890
$ws.Range("M99").Value = 608

$ws.Range("H107").Value = 1770.1428
$ws.Range("I107").Value = 1765.1666
$ws.Range("J107").Value = 1800
$ws.Range("K107").Value = 1765.1666
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 154.8334
$ws.Range("N107").Value = -5640

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 1650
$ws.Range("I48").Value = 300
$ws.Range("J48").Value = 3000
$ws.Range("K48").Value = 900
$ws.Range("L48").Value = 9000
$ws.Range("M48").Value = -650
$ws.Range("N48").Value = -9500

$ws.Range("H131").Value = 3379.75
$ws.Range("J131").Value = 2394
$ws.Range("L131").Value = 7182
$ws.Range("N131").Value = -17262

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 55557136
$ws.Range("I97").Value = 83334696
$ws.Range("J97").Value = 2026.8334
$ws.Range("K97").Value = 83334696
$ws.Range("L97").Value = 2026.8334
$ws.Range("M97").Value = -83334200
$ws.Range("N97").Value = -3018.8334

$ws.Range("H122").Value = 2824.4644
$ws.Range("I122").Value = 2000.909
$ws.Range("J122").Value = 3357.353
$ws.Range("K122").Value = 6002.727000000001
$ws.Range("L122").Value = 10072.059
$ws.Range("M122").Value = -3552.727000000001
$ws.Range("N122").Value = -14972.059

$ws.Range("H132").Value = 1609
$ws.Range("I132").Value = 1227.1538
$ws.Range("J132").Value = 2318.1428
$ws.Range("K132").Value = 3681.4614
$ws.Range("L132").Value = 6954.428400000001
$ws.Range("M132").Value = -1151.4614
$ws.Range("N132").Value = -12014.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2383.5293
$ws.Range("I40").Value = 2452.8333
$ws.Range("J40").Value = 2217.2
$ws.Range("K40").Value = 2452.8333
$ws.Range("L40").Value = 2217.2
$ws.Range("M40").Value = -2316.8333
$ws.Range("N40").Value = -2489.2

$ws.Range("H61").Value = 1741
$ws.Range("I61").Value = 1356
$ws.Range("J61").Value = 2280
$ws.Range("K61").Value = 1356
$ws.Range("L61").Value = 2280
$ws.Range("M61").Value = -1154
$ws.Range("N61").Value = -2684

$ws.Range("H93").Value = 1730.6
$ws.Range("I93").Value = 1801.2
$ws.Range("J93").Value = 1660
$ws.Range("K93").Value = 1801.2
$ws.Range("L93").Value = 1660
$ws.Range("M93").Value = -553.2
$ws.Range("N93").Value = -4156

$ws.Range("H100").Value = 2564.7144
$ws.Range("I100").Value = 1875.75
$ws.Range("J100").Value = 3483.3333
$ws.Range("K100").Value = 1875.75
$ws.Range("L100").Value = 3483.3333
$ws.Range("M100").Value = -1334.75
$ws.Range("N100").Value = -4565.3333

$ws.Range("H113").Value = 1741
$ws.Range("I113").Value = 1356
$ws.Range("J113").Value = 2280
$ws.Range("K113").Value = 1356
$ws.Range("L113").Value = 2280
$ws.Range("M113").Value = 814
$ws.Range("N113").Value = -6620

$ws.Range("H122").Value = 2687
$ws.Range("I122").Value = 2440
$ws.Range("J122").Value = 3366.25
$ws.Range("K122").Value = 7320
$ws.Range("L122").Value = 10098.75
$ws.Range("M122").Value = -4870
$ws.Range("N122").Value = -14998.75

$ws.Range("H132").Value = 5640.7812
$ws.Range("I132").Value = 8672.5
$ws.Range("J132").Value = 3282.7778
$ws.Range("K132").Value = 26017.5
$ws.Range("L132").Value = 9848.3334
$ws.Range("M132").Value = -23487.5
$ws.Range("N132").Value = -14908.3334

$ws.Range("H136").Value = 1837.6154
$ws.Range("I136").Value = 1087.8235
$ws.Range("K136").Value = 3263.4705
$ws.Range("M136").Value = -713.4704999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 424.5
$ws.Range("I107").Value = 399.33334
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 1198.00002
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 721.9999800000001
$ws.Range("N107").Value = -5340

$ws.Range("H122").Value = 1851.9756
$ws.Range("I122").Value = 1232.7693
$ws.Range("J122").Value = 2925.2666
$ws.Range("K122").Value = 3698.3079
$ws.Range("L122").Value = 8775.799800000001
$ws.Range("M122").Value = -1248.3079
$ws.Range("N122").Value = -13675.7998

$ws.Range("H136").Value = 13223.647
$ws.Range("I136").Value = 16784.385
$ws.Range("J136").Value = 1651.25
$ws.Range("K136").Value = 50353.155
$ws.Range("L136").Value = 4953.75
$ws.Range("M136").Value = -47803.155
$ws.Range("N136").Value = -10053.75
